# Generate Report for Handback
# Update the timestamp values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: G2 is "Latest HO Xliff Generate Date" for the first file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 01:12:49"

# "zh-cn" sheet: H2 is "Correspond Handoff Datetime", K2 is "Correspond Handback DateTime"
# for the first file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 01:12:44"
$wsZhCn.Range("K2").Value = "2016-09-06 01:13:02"

# "de-de" sheet: H2 is "Correspond Handoff Datetime" (shares the Overview!G2 timestamp),
# K2 is "Correspond Handback DateTime" for the first file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 01:12:49"
$wsDeDe.Range("K2").Value = "2016-09-06 01:13:15"
